# 🔄 MAJ automatique BRVM via GitHub Actions
#
# This script refreshes the "Recommandations" and "Top_YTD" worksheets with
# the latest BRVM market data (jours en hausse/baisse, variations, YTD, etc.)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Recommandations": columns A..G, data rows 2..51
# Columns: A=Titre, B=Jours en Hausse, C=Jours en Baisse,
#          D=Variation Totale (%), E=Derniere Variation (%),
#          F=Recommandation, G=Strategie
# ------------------------------------------------------------------
$wsReco = $wb.Worksheets.Item("Recommandations")

$recoData = @(
    @(2, 'BRVM - SERVICES PUBLICS', 0, 10, 4154.56, 109.02, '🟡 Observer', '➖ Neutre'),
    @(3, 'SAFCA CI', 0, 5, 3465, 685, '🟡 Observer', '➖ Neutre'),
    @(4, 'CFAO MOTORS CI', 0, 5, 3400, 680, '🟡 Observer', '➖ Neutre'),
    @(5, 'BRVM - AUTRES SECTEURS', 0, 5, 3318.3, 663.3200000000001, '🟡 Observer', '➖ Neutre'),
    @(6, 'NEI-CEDA CI', 0, 5, 2960, 600, '🟡 Observer', '➖ Neutre'),
    @(7, 'SETAO CI', 0, 5, 2875, 595, '🟡 Observer', '➖ Neutre'),
    @(8, 'UNIWAX CI', 0, 5, 2825, 565, '🟡 Observer', '➖ Neutre'),
    @(9, 'AIR LIQUIDE CI', 0, 5, 2610, 520, '🟡 Observer', '➖ Neutre'),
    @(10, 'BRVM - DISTRIBUTION', 0, 5, 1828.82, 356.71, '🟡 Observer', '➖ Neutre'),
    @(11, 'BRVM - TRANSPORT', 0, 5, 1748.91, 348.8, '🟡 Observer', '➖ Neutre'),
    @(12, 'BRVM - AGRICULTURE', 0, 5, 1528.52, 306.86, '🟡 Observer', '➖ Neutre'),
    @(13, 'BRVM - INDUSTRIE', 0, 5, 1045.86, 214.12, '🟡 Observer', '➖ Neutre'),
    @(14, 'SUCRIVOIRE', 0, 1, 995, 995, '🟡 Observer', '➖ Neutre'),
    @(15, 'BRVM - CONSOMMATION DE BASE', 0, 5, 891.13, 181.66, '🟡 Observer', '➖ Neutre'),
    @(16, 'BRVM-PRINCIPAL', 0, 5, 885.89, 177.9, '🟡 Observer', '➖ Neutre'),
    @(17, 'BRVM - INDUSTRIELS', 0, 5, 657.51, 131.75, '🟡 Observer', '➖ Neutre'),
    @(18, 'BRVM-PRESTIGE', 0, 5, 651.3200000000001, 129.18, '🟡 Observer', '➖ Neutre'),
    @(19, 'BRVM - FINANCES', 0, 5, 611.16, 121.39, '🟡 Observer', '➖ Neutre'),
    @(20, 'BRVM - SERVICES FINANCIERS', 0, 5, 600.64, 119.3, '🟡 Observer', '➖ Neutre'),
    @(21, 'BRVM - ENERGIE', 0, 5, 536.6799999999999, 105.76, '🟡 Observer', '➖ Neutre'),
    @(22, 'BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 5, 534.34, 104.95, '🟡 Observer', '➖ Neutre'),
    @(23, 'BRVM - TELECOMMUNICATIONS', 0, 5, 470.7, 94.25, '🟡 Observer', '➖ Neutre'),
    @(24, 'UNILEVER CI (UNLC)', 4, 0, 29.95, 7.5, '🟢 Achat', '✅ Renforcer'),
    @(25, 'BANK OF AFRICA SENEGAL (BOAS)', 2, 0, 8.289999999999999, 2.3, '🟡 Observer', '➖ Neutre'),
    @(26, 'SETAO CI (STAC)', 1, 0, 7.27, 7.27, '🟡 Observer', '➖ Neutre'),
    @(27, 'SAFCA CI (SAFC)', 1, 0, 6.92, 6.92, '🟡 Observer', '➖ Neutre'),
    @(28, 'SOLIBRA CI (SLBC)', 1, 0, 5.72, 5.72, '🟡 Observer', '➖ Neutre'),
    @(29, 'BICI CI (BICC)', 1, 0, 5.7, 5.7, '🟡 Observer', '➖ Neutre'),
    @(30, 'NSIA BANQUE COTE D''IVOIRE (NSBC)', 2, 0, 3.97, 1.54, '🟡 Observer', '➖ Neutre'),
    @(31, 'SMB CI (SMBC)', 2, 1, 3.83, 3.92, '🟡 Observer', '👀 À surveiller'),
    @(32, 'ORANGE COTE D''IVOIRE (ORAC)', 1, 0, 3.57, 3.57, '🟡 Observer', '➖ Neutre'),
    @(33, 'BERNABE CI (BNBC)', 3, 2, 2.45, 3.46, '🟡 Observer', '➖ Neutre'),
    @(34, 'UNIWAX CI (UNXC)', 1, 1, 2.12, 6.67, '🟡 Observer', '👀 À surveiller'),
    @(35, 'ONATEL BF (ONTBF)', 1, 1, 2.01, 5.18, '🟡 Observer', '👀 À surveiller'),
    @(36, 'SAPH CI (SPHC)', 1, 1, 1.13, -2.13, '🟡 Observer', '👀 À surveiller'),
    @(37, 'ECOBANK TRANS. INCORP. TG (ETIT)', 1, 1, 0.37, -5.88, '🟡 Observer', '👀 À surveiller'),
    @(38, 'TOTAL', 0, 5, 0, 0, '🟡 Observer', '➖ Neutre'),
    @(39, 'SODE CI (SDCC)', 1, 2, -1.66, 7.26, '🟡 Observer', '👀 À surveiller'),
    @(40, 'AIR LIQUIDE CI (SIVC)', 1, 2, -1.8, -2.86, '🟡 Observer', '👀 À surveiller'),
    @(41, 'ORAGROUP TOGO (ORGT)', 0, 1, -2.13, -2.13, '🟡 Observer', '➖ Neutre'),
    @(42, 'VIVO ENERGY CI (SHEC)', 0, 1, -2.4, -2.4, '🟡 Observer', '➖ Neutre'),
    @(43, 'TOTALENERGIES MARKETING CI (TTLC)', 0, 1, -2.44, -2.44, '🟡 Observer', '➖ Neutre'),
    @(44, 'TOTALENERGIES MARKETING SN (TTLS)', 0, 1, -2.57, -2.57, '🟡 Observer', '➖ Neutre'),
    @(45, 'CIE CI (CIEC)', 1, 2, -2.92, 7.33, '🟡 Observer', '👀 À surveiller'),
    @(46, 'CFAO MOTORS CI (CFAC)', 0, 1, -2.99, -2.99, '🟡 Observer', '➖ Neutre'),
    @(47, 'SOGB CI (SOGC)', 0, 1, -3.04, -3.04, '🟡 Observer', '➖ Neutre'),
    @(48, 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)', 0, 2, -3.78, -1.86, '🟡 Observer', '➖ Neutre'),
    @(49, 'SICOR CI (SICC)', 0, 1, -5.71, -5.71, '🟡 Observer', '➖ Neutre'),
    @(50, 'ECOBANK COTE D''''IVOIRE (ECOC)', 0, 1, -5.93, -5.93, '🟡 Observer', '➖ Neutre'),
    @(51, 'TRACTAFRIC MOTORS CI (PRSC)', 0, 2, -8.94, -5.9, '🟡 Observer', '➖ Neutre')
)

foreach ($row in $recoData) {
    $r = $row[0]
    $wsReco.Cells.Item($r, 1).Value = $row[1]
    $wsReco.Cells.Item($r, 2).Value = $row[2]
    $wsReco.Cells.Item($r, 3).Value = $row[3]
    $wsReco.Cells.Item($r, 4).Value = $row[4]
    $wsReco.Cells.Item($r, 5).Value = $row[5]
    $wsReco.Cells.Item($r, 6).Value = $row[6]
    $wsReco.Cells.Item($r, 7).Value = $row[7]
}

# ------------------------------------------------------------------
# Sheet "Top_YTD": columns A..B, data rows 2..11
# Columns: A=Titre, B=Progression YTD (%)
# Only column B (values) changes for this refresh.
# ------------------------------------------------------------------
$wsTop = $wb.Worksheets.Item("Top_YTD")

$topData = @(
    @(2, 147324810.79),
    @(3, 3135502.43),
    @(4, 2886716.93),
    @(5, 2597058.77),
    @(6, 1586598.47),
    @(7, 1396349.6),
    @(8, 1299875.67),
    @(9, 930748.62),
    @(10, 219010.66),
    @(11, 183955.14)
)

foreach ($row in $topData) {
    $r = $row[0]
    $wsTop.Cells.Item($r, 2).Value = $row[1]
}
